$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: column B "A" -> "C", column C "C" -> "FFR"
$ws.Range("B1").Value = "C"
$ws.Range("C1").Value = "FFR"

# Row labels: "A Lag" -> "C Lag", "C Lag" -> "FFR Lag"
$ws.Range("A2").Value = "C Lag"
$ws.Range("A3").Value = "FFR Lag"

# Data cells B2:D4 must hold TEXT values, even the ones that look like
# plain numbers ("3.79", "-6.09", "-0.01"). A direct Range.Value = "3.79"
# assignment gets auto-coerced to a numeric cell by Excel, so for those
# purely-numeric-looking entries we instead write a =TEXT() formula and
# immediately convert it to a static value via copy / paste-special
# (values only) - this yields a genuine text cell without ever touching
# NumberFormat (which would otherwise leave a permanent, if unused, extra
# style registered in the workbook).

# Row 2 ("C Lag"): -0.46***, 3.79, -6.09
$ws.Range("B2").Value = "-0.46***"
$ws.Range("C2").Formula = '=TEXT(3.79,"0.00")'
$ws.Range("D2").Formula = '=TEXT(-6.09,"0.00")'
$ws.Range("C2:D2").Copy()
$ws.Range("C2:D2").PasteSpecial(-4163) # xlPasteValues

# Row 3 ("FFR Lag"): -0.01, 1.6***, 0.5***
$ws.Range("B3").Formula = '=TEXT(-0.01,"0.00")'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163) # xlPasteValues
$ws.Range("C3").Value = "1.6***"
$ws.Range("D3").Value = "0.5***"

# Row 4 ("LF Lag"): 0.04*, 3.53*, 0.54*
$ws.Range("B4").Value = "0.04*"
$ws.Range("C4").Value = "3.53*"
$ws.Range("D4").Value = "0.54*"

$excel.CutCopyMode = $false
